# Update the confidence-interval text values in columns M (CI(amplitude))
# and U (CI(mesor)) for rows 2-10 with re-computed figures from a re-run
# of the CircadiPy / CircaDB cosinor analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "M2"  = "[49.86848824929016, 50.09098322812482]"
    "U2"  = "[49.9630732407347, 50.12345618101146]"
    "M3"  = "[49.91596268719528, 50.20114184090381]"
    "U3"  = "[49.91957618909732, 50.095208606307345]"
    "M4"  = "[49.8688696571068, 50.17773194296102]"
    "U4"  = "[49.82872672181338, 49.99837747241425]"
    "M5"  = "[49.88708651385721, 50.1661080751208]"
    "U5"  = "[49.894313088102066, 50.037645845375465]"
    "M6"  = "[49.80595047620726, 50.146339927165755]"
    "U6"  = "[49.878854617897694, 50.05462834506002]"
    "M7"  = "[49.85349469012344, 50.169234847215705]"
    "U7"  = "[49.93691333315288, 50.108111745779944]"
    "M8"  = "[49.80520944595636, 50.09202248830258]"
    "U8"  = "[49.929931324276644, 50.11740060402587]"
    "M9"  = "[49.758185987673386, 50.077329237853206]"
    "U9"  = "[49.925853921652894, 50.106676693751716]"
    "M10" = "[49.730321754939936, 50.062306237414504]"
    "U10" = "[49.88548510526114, 50.05565384015115]"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
